$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 2: Taxonsorteringsordning changes
$ws.Range("B2").Value = 79244

# Row 3: Taxonsorteringsordning changes
$ws.Range("B3").Value = 79244

# Row 4: Id, Taxonsorteringsordning, Ost, Nord and Publik kommentar change
$ws.Range("A4").Value = 131082789
$ws.Range("B4").Value = 79244
$ws.Range("Q4").Value = 459958
$ws.Range("R4").Value = 7046508
$ws.Range("AC4").Value = "På en stående död gran."

# Row 5: Id, Taxonsorteringsordning, Ost, Nord and Publik kommentar change
$ws.Range("A5").Value = 131082790
$ws.Range("B5").Value = 79244
$ws.Range("Q5").Value = 459913
$ws.Range("R5").Value = 7046493
$ws.Range("AC5").Value = "På flera stående döda granar."
